$wb = $excel.ActiveWorkbook

# --- "VerifyStartEndDateValidation" sheet (currently the active tab) ---
# Clear the stale FAIL result and its comment message, then move the
# selection off this sheet (it stops being the active tab).
$wsVal = $wb.Worksheets.Item("VerifyStartEndDateValidation")
$wsVal.Range("K2").Value = ""
$wsVal.Range("L2").Value = ""
$wsVal.Range("I2").Select() | Out-Null

# --- "Test Cases" summary sheet ---
# Clear the rolled-up Result value (FAIL) for the
# VerifyStartEndDateValidation row, then make this sheet the active tab.
$wsTC = $wb.Worksheets.Item("Test Cases")
$wsTC.Range("F2").Value = ""
$wsTC.Activate() | Out-Null
$wsTC.Range("D2").Select() | Out-Null
